$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data changes ---

# Row 8: fill in G8 and H8 with 5 (J8 formula auto-recalculates to 30)
$ws.Range("G8").Value = 5
$ws.Range("H8").Value = 5

# Row 11: fill in G11 and H11 with 5 (J11 formula auto-recalculates to 30)
$ws.Range("G11").Value = 5
$ws.Range("H11").Value = 5

# Row 16: C16:F16 change from 0 to 5, and restyle to match "s=2" (green fill) cells
$ws.Range("C16:F16").Value = 5
$ws.Range("C16:F16").Style = $ws.Range("G16").Style

# --- View changes ---

# Zoom to 190%
$excel.ActiveWindow.Zoom = 190

# Move frozen pane's top-left visible cell to D4, and set active selection to I8
$ws.Range("I8").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 4
